$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data refresh of the cryptos table (price + 1h volume change columns),
# including a few rows whose coin ordering swapped with their neighbour.

# Row 2
$ws.Range("D2").Value = "'26.153.03"
$ws.Range("E2").Value = "'  -4.15%  "
# Row 3
$ws.Range("D3").Value = "'1.650.69"
$ws.Range("E3").Value = "'  -3.51%  "
# Row 4
$ws.Range("E4").Value = "'  +0.17%  "
# Row 5
$ws.Range("D5").Value = "'215.54"
$ws.Range("E5").Value = "'  -4.16%  "
# Row 6
$ws.Range("D6").Value = "'0.5125"
$ws.Range("E6").Value = "'  -3.13%  "
# Row 7
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "'  +0.23%  "
# Row 8
$ws.Range("D8").Value = "'0.2591"
$ws.Range("E8").Value = "'  -2.60%  "
# Row 9
$ws.Range("D9").Value = "'0.06431"
$ws.Range("E9").Value = "'  -3.86%  "
# Row 10
$ws.Range("D10").Value = "'19.92"
$ws.Range("E10").Value = "'  -4.51%  "
# Row 11
$ws.Range("D11").Value = "'0.07793"
$ws.Range("E11").Value = "'  +1.35%  "
# Row 12
$ws.Range("D12").Value = "'1.650.54"
$ws.Range("E12").Value = "'  -3.51%  "
# Row 13
$ws.Range("E13").Value = "'  -4.73%  "
# Row 14
$ws.Range("D14").Value = "'1.877.10"
$ws.Range("E14").Value = "'  -3.54%  "
# Row 15
$ws.Range("D15").Value = "'0.5512"
$ws.Range("E15").Value = "'  -5.87%  "
# Row 16
$ws.Range("D16").Value = "'0.0₅8003"
$ws.Range("E16").Value = "'  -2.63%  "
# Row 17
$ws.Range("E17").Value = "'  -5.66%  "
# Row 18
$ws.Range("D18").Value = "'26.152.09"
$ws.Range("E18").Value = "'  -4.28%  "
# Row 19
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "'  +0.16%  "
# Row 20
$ws.Range("D20").Value = "'210.28"
$ws.Range("E20").Value = "'  -5.18%  "
# Row 21
$ws.Range("D21").Value = "'4.399"
$ws.Range("E21").Value = "'  -5.37%  "
# Row 22
$ws.Range("E22").Value = "'  -4.03%  "
# Row 23
$ws.Range("D23").Value = "'6.037"
$ws.Range("E23").Value = "'  +0.12%  "
# Row 24
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "'  +0.14%  "
# Row 25
$ws.Range("D25").Value = "'144.01"
$ws.Range("E25").Value = "'  -0.60%  "
# Row 26
$ws.Range("E26").Value = "'  +3.11%  "
# Row 27
$ws.Range("D27").Value = "'0.1175"
$ws.Range("E27").Value = "'  -2.74%  "
# Row 28
$ws.Range("D28").Value = "'6.963"
$ws.Range("E28").Value = "'  -3.75%  "
# Row 29
$ws.Range("D29").Value = "'15.80"
$ws.Range("E29").Value = "'  -2.78%  "
# Row 30
$ws.Range("D30").Value = "'0.05132"
$ws.Range("E30").Value = "'  -3.85%  "
# Row 31
$ws.Range("D31").Value = "'1.241"
$ws.Range("E31").Value = "'  -4.01%  "
# Row 32
$ws.Range("D32").Value = "'3.350"
$ws.Range("E32").Value = "'  -3.50%  "
# Row 33
$ws.Range("D33").Value = "'3.213"
$ws.Range("E33").Value = "'  -6.20%  "
# Row 34
$ws.Range("D34").Value = "'1.558"
$ws.Range("E34").Value = "'  -4.37%  "
# Row 35
$ws.Range("D35").Value = "'2.737"
$ws.Range("E35").Value = "'  -4.65%  "
# Row 36
$ws.Range("B36").Value = "'ARBITRUM"
$ws.Range("C36").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9240"
$ws.Range("E36").Value = "'  -3.35%  "
# Row 37
$ws.Range("B37").Value = "'HuobiToken"
$ws.Range("C37").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.352"
$ws.Range("E37").Value = "'  -1.65%  "
# Row 38
$ws.Range("B38").Value = "'ImmutableX"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.5704"
$ws.Range("E38").Value = "'  -2.67%  "
# Row 39
$ws.Range("B39").Value = "'Maker"
$ws.Range("C39").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "'1.162.44"
$ws.Range("E39").Value = "'  +1.45%  "
# Row 40
$ws.Range("E40").Value = "'  -3.28%  "
# Row 41
$ws.Range("B41").Value = "'mCoin"
$ws.Range("C41").Value = "'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D41").Value = "'2.560"
$ws.Range("E41").Value = "'  -0.12%  "
# Row 42
$ws.Range("B42").Value = "'PaxDollar"
$ws.Range("C42").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.005"
$ws.Range("E42").Value = "'  +0.11%  "
# Row 43
$ws.Range("D43").Value = "'5.652"
$ws.Range("E43").Value = "'  -2.35%  "
# Row 44
$ws.Range("D44").Value = "'0.8226"
$ws.Range("E44").Value = "'  -1.88%  "
# Row 45
$ws.Range("D45").Value = "'100.09"
$ws.Range("E45").Value = "'  -0.81%  "
# Row 46
$ws.Range("D46").Value = "'1.788.43"
$ws.Range("E46").Value = "'  -3.49%  "
# Row 47
$ws.Range("D47").Value = "'0.0₈116"
$ws.Range("E47").Value = "'  +3.94%  "
# Row 48
$ws.Range("E48").Value = "'  -0.44%  "
# Row 49
$ws.Range("E49").Value = "'  -3.72%  "
# Row 50
$ws.Range("D50").Value = "'1.005"
$ws.Range("E50").Value = "'  +0.23%  "
# Row 51
$ws.Range("D51").Value = "'7.831"
$ws.Range("E51").Value = "'  -3.18%  "
